$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift rows 7..11 down to 8..12 (bottom-up to avoid clobbering sources) ---
# Each row is copied in two passes (formats, then values) because this engine's
# PasteSpecial(xlPasteValues) is a no-op when the source range is entirely blank,
# which would otherwise leave stale values behind; ClearContents() beforehand
# guarantees the destination is actually blanked in that case.
for ($r = 11; $r -ge 7; $r--) {
    $dst = $r + 1
    $srcRange = $ws.Range("A" + $r + ":P" + $r)
    $dstRange = $ws.Range("A" + $dst + ":P" + $dst)

    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)   # xlPasteFormats
    $dstRange.ClearContents()

    $srcRange.Copy()
    $dstRange.PasteSpecial(-4163)   # xlPasteValues

    $dstRange.RowHeight = $srcRange.RowHeight
}

# At this point row 6 still holds the original "SONAR_ON" row untouched, and row 7
# still holds a duplicate of the original "GREEN" row (nothing wrote into row 6 or
# row 7 as a destination above). Re-purpose row 6 as "SONAR FAR" and row 7 as the
# brand-new "SONAR_NEAR" row.

# --- Re-skin row 6 using row 8's pattern (same visual family), keep the checkmark in column K ---
$ws.Range("A8:P8").Copy()
$ws.Range("A6:P6").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B8").Copy()
$ws.Range("I6").PasteSpecial(-4122)      # restore I6 to the plain (unchecked) style

$ws.Range("F11").Copy()
$ws.Range("K6").PasteSpecial(-4122)      # checked style for column K

$ws.Range("A6").Value = "SONAR FAR"
$ws.Range("K6").Value = "C"

# --- Build row 7 as the new "SONAR_NEAR" row, using row 9's pattern, checkmark in column J ---
$ws.Range("A9:P9").Copy()
$ws.Range("A7:P7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A7:P7").ClearContents()

$ws.Range("F11").Copy()
$ws.Range("J7").PasteSpecial(-4122)      # checked style for column J

$ws.Range("A7").Value = "SONAR_NEAR"
$ws.Range("J7").Value = "C"

# --- Update the header selection to match the target state ---
$ws.Range("A7").Select()

Write-Output "done"
